$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Add new row 4 with the new user's data
$ws.Range("A4").Value = '"5f05af6804dfeb226c3c42df"'
$ws.Range("B4").Value = "Efehi"
$ws.Range("C4").Value = "efehi@gmail.com"
$ws.Range("D4").Value = 918876543210
